# Atualização de TOKEN JWT
# Updates the "FOLHA 1" comparison table (rows 17-31) to reflect the new
# revision comparison data: a new item is inserted at row 17
# (10027489 / tê para tubo), a new item is inserted at row 22
# (10334519 / válvula gaveta), all following rows shift down by one,
# and three additional rows (29-31) are populated with items that were
# previously empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FOLHA 1")

# Table of new values for columns E (code), L (description), AV (qtd rev2),
# BB (qtd rev4), BH (diferença), BM (status) for rows 17 through 31.
$rows = @(
    @{ Row=17; A=11; E="10027489"; L="tê para tubo";                    AV=0;   BB=2;                BH=2;                   BM="ADICIONADO" },
    @{ Row=18; A=12; E="10038562"; L="colar para tubo";                 AV=1;   BB=2;                BH=1;                   BM="ALTERADO" },
    @{ Row=19; A=13; E="10217343"; L="Flange pescoo";                   AV=4;   BB=5;                BH=1;                   BM="ALTERADO" },
    @{ Row=20; A=14; E="10290551"; L="Flange pescoo";                   AV=2;   BB=5;                BH=3;                   BM="ALTERADO" },
    @{ Row=21; A=15; E="10312742"; L="junta circular para flange";     AV=1;   BB=2;                BH=1;                   BM="ALTERADO" },
    @{ Row=22; A=16; E="10334519"; L="válvula gaveta";                  AV=0;   BB=4;                BH=4;                   BM="ADICIONADO" },
    @{ Row=23; A=17; E="10360352"; L="codigo gp";                       AV=3;   BB=6;                BH=3;                   BM="ALTERADO" },
    @{ Row=24; A=18; E="10390568"; L="tubo condução de aço carbono";   AV=4.9; BB=4.8999999999999995; BH=-8.881784197001252e-16; BM="ALTERADO" },
    @{ Row=25; A=19; E="10390629"; L="tubo condução de aço carbono";   AV=1.6; BB=2.2;              BH=0.6000000000000001; BM="ALTERADO" },
    @{ Row=26; A=20; E="10515625"; L="junta circular para flange";     AV=2;   BB=6;                BH=4;                   BM="ALTERADO" },
    @{ Row=27; A=21; E="10559349"; L="parafuso estojo";                 AV=16;  BB=24;               BH=8;                   BM="ALTERADO" },
    @{ Row=28; A=22; E="11389612"; L="junta circular para flange";     AV=2;   BB=5;                BH=3;                   BM="ALTERADO" },
    @{ Row=29; A=23; E="11389613"; L="junta circular para flange";     AV=4;   BB=0;                BH=-4;                  BM="ALTERADO" },
    @{ Row=30; A=24; E="11389614"; L="junta circular para flange";     AV=6;   BB=0;                BH=-6;                  BM="ALTERADO" },
    @{ Row=31; A=25; E="11399557"; L="tubo condução de aço carbono";   AV=0;   BB=0.6;              BH=0.6;                 BM="ADICIONADO" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A        # column A
    $ws.Cells.Item($n, 5).Value = $r.E        # column E
    $ws.Cells.Item($n, 12).Value = $r.L       # column L
    $ws.Cells.Item($n, 48).Value = $r.AV      # column AV
    $ws.Cells.Item($n, 54).Value = $r.BB      # column BB
    $ws.Cells.Item($n, 60).Value = $r.BH      # column BH
    $ws.Cells.Item($n, 65).Value = $r.BM      # column BM
}

# Row 11: only the "QTD REV(4)" (BB) and "DIFERENÇA" (BH) values changed.
$ws.Cells.Item(11, 54).Value = 2
$ws.Cells.Item(11, 60).Value = 1
